$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H19").Value = 361
$ws.Range("I19").Value = 487.5
$ws.Range("J19").Value = 276.66666
$ws.Range("K19").Value = 487.5
$ws.Range("L19").Value = 276.66666
$ws.Range("M19").Value = -312.5
$ws.Range("N19").Value = -626.66666

# ALC row 51
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H51").Value = 4420.4
$ws.Range("J51").Value = 4950.5
$ws.Range("L51").Value = 4950.5
$ws.Range("N51").Value = -5918.5

# ALC row 125
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H125").Value = 1332.4
$ws.Range("I125").Value = 1061.6
$ws.Range("J125").Value = 1467.8
$ws.Range("K125").Value = 9554.4
$ws.Range("L125").Value = 13210.2
$ws.Range("M125").Value = -7094.4
$ws.Range("N125").Value = -18130.2

# ALC row 132
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H132").Value = 5964258
$ws.Range("I132").Value = 11383501
$ws.Range("J132").Value = 3090.5
$ws.Range("K132").Value = 34150503
$ws.Range("L132").Value = 9271.5
$ws.Range("M132").Value = -34147973
$ws.Range("N132").Value = -14331.5

# ALC row 137
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H137").Value = 4929.3335
$ws.Range("I137").Value = 4956.5625
$ws.Range("J137").Value = 4874.875
$ws.Range("K137").Value = 14869.6875
$ws.Range("L137").Value = 14624.625
$ws.Range("M137").Value = -12319.6875
$ws.Range("N137").Value = -19724.625

# ARM row 74
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H74").Value = 18144.857
$ws.Range("I74").Value = 50450
$ws.Range("K74").Value = 50450
$ws.Range("M74").Value = -49576

# ARM row 77
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H77").Value = 18144.857
$ws.Range("I77").Value = 50450
$ws.Range("K77").Value = 252250
$ws.Range("M77").Value = -247882

# ARM row 80
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H80").Value = 30832.5
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 40110
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 40110
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -42106

# ARM row 83
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H83").Value = 30832.5
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 40110
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 120330
$ws.Range("M83").Value = -4008
$ws.Range("N83").Value = -130314

# ARM row 132
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H132").Value = 140312.25
$ws.Range("I132").Value = 1000000
$ws.Range("J132").Value = 17499.715
$ws.Range("K132").Value = 3000000
$ws.Range("L132").Value = 52499.145
$ws.Range("M132").Value = -2997470
$ws.Range("N132").Value = -57559.145

# BSM row 134
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H134").Value = 7710.294
$ws.Range("I134").Value = 8791.071
$ws.Range("J134").Value = 2666.6667
$ws.Range("K134").Value = 26373.213
$ws.Range("L134").Value = 8000.000100000001
$ws.Range("M134").Value = -23838.213
$ws.Range("N134").Value = -13070.0001

# CRP row 31
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 4613.3237
$ws.Range("I31").Value = 6102.8945
$ws.Range("K31").Value = 6102.8945
$ws.Range("M31").Value = -5807.8945

# CRP row 34
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H34").Value = 4613.3237
$ws.Range("I34").Value = 6102.8945
$ws.Range("K34").Value = 6102.8945
$ws.Range("M34").Value = -5900.8945

# CRP row 52
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H52").Value = 49660
$ws.Range("J52").Value = 49660
$ws.Range("L52").Value = 49660
$ws.Range("N52").Value = -50248

# CRP row 58
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H58").Value = 985.75
$ws.Range("I58").Value = 985.75
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 985.75
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -782.75

# CRP row 99
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H99").Value = 1436.75
$ws.Range("I99").Value = 1385.5555
$ws.Range("J99").Value = 1502.5714
$ws.Range("K99").Value = 1385.5555
$ws.Range("L99").Value = 1502.5714
$ws.Range("M99").Value = 112.4445000000001
$ws.Range("N99").Value = -4498.5714

# CRP row 126
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H126").Value = 1436.75
$ws.Range("I126").Value = 1385.5555
$ws.Range("J126").Value = 1502.5714
$ws.Range("K126").Value = 4156.666499999999
$ws.Range("L126").Value = 4507.7142
$ws.Range("M126").Value = -1686.666499999999
$ws.Range("N126").Value = -9447.7142

# CRP row 132
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H132").Value = 11268.546
$ws.Range("I132").Value = 7744.375
$ws.Range("K132").Value = 23233.125
$ws.Range("M132").Value = -20703.125

# CRP row 134
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H134").Value = 11847.5
$ws.Range("I134").Value = 10330.546
$ws.Range("J134").Value = 13701.556
$ws.Range("K134").Value = 30991.638
$ws.Range("L134").Value = 41104.66800000001
$ws.Range("M134").Value = -28456.638
$ws.Range("N134").Value = -46174.66800000001

# CRP row 136
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H136").Value = 985.75
$ws.Range("I136").Value = 985.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2957.25
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -407.25

# CUL row 5
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H5").Value = 714.25714
$ws.Range("I5").Value = 438.08334
$ws.Range("K5").Value = 1314.25002
$ws.Range("M5").Value = -1202.25002

# CUL row 12
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H12").Value = 103.63158
$ws.Range("I12").Value = 92.5
$ws.Range("J12").Value = 116
$ws.Range("K12").Value = 277.5
$ws.Range("L12").Value = 348
$ws.Range("M12").Value = -104.5
$ws.Range("N12").Value = -694

# CUL row 113
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H113").Value = 944.9167
$ws.Range("I113").Value = 1375.0588
$ws.Range("J113").Value = 560.0526
$ws.Range("K113").Value = 4125.1764
$ws.Range("L113").Value = 1680.1578
$ws.Range("M113").Value = -1955.1764
$ws.Range("N113").Value = -6020.1578

# CUL row 118
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H118").Value = 1791.0714
$ws.Range("I118").Value = 382.6
$ws.Range("K118").Value = 1147.8
$ws.Range("M118").Value = 95.19999999999982

# CUL row 135
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H135").Value = 714.25714
$ws.Range("I135").Value = 438.08334
$ws.Range("K135").Value = 3942.75006
$ws.Range("M135").Value = -1407.75006

# GSM row 22
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H22").Value = 32504
$ws.Range("I22").Value = 5008
$ws.Range("J22").Value = 60000
$ws.Range("K22").Value = 5008
$ws.Range("L22").Value = 60000
$ws.Range("M22").Value = -4479
$ws.Range("N22").Value = -61058

# GSM row 137
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H137").Value = 70000
$ws.Range("J137").Value = 70000
$ws.Range("L137").Value = 70000
$ws.Range("N137").Value = -80200

# LTW row 7
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H7").Value = 2751.353
$ws.Range("I7").Value = 2246.9167
$ws.Range("J7").Value = 3962
$ws.Range("K7").Value = 2246.9167
$ws.Range("L7").Value = 3962
$ws.Range("M7").Value = -2134.9167
$ws.Range("N7").Value = -4186

# LTW row 16
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H16").Value = 825.0909
$ws.Range("J16").Value = 746.8
$ws.Range("L16").Value = 746.8
$ws.Range("N16").Value = -1086.8

# LTW row 126
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H126").Value = 2751.353
$ws.Range("I126").Value = 2246.9167
$ws.Range("J126").Value = 3962
$ws.Range("K126").Value = 6740.750100000001
$ws.Range("L126").Value = 11886
$ws.Range("M126").Value = -4270.750100000001
$ws.Range("N126").Value = -16826

# WVR row 122
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H122").Value = 3354508
$ws.Range("I122").Value = 1954346
$ws.Range("J122").Value = 5954809
$ws.Range("K122").Value = 5863038
$ws.Range("L122").Value = 17864427
$ws.Range("M122").Value = -5860588
$ws.Range("N122").Value = -17869327

# WVR row 126
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H126").Value = 1272416
$ws.Range("I126").Value = 1211626.1
$ws.Range("J126").Value = 1359801.5
$ws.Range("K126").Value = 3634878.3
$ws.Range("L126").Value = 4079404.5
$ws.Range("M126").Value = -3632408.3
$ws.Range("N126").Value = -4084344.5

# WVR row 132
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H132").Value = 9289
$ws.Range("I132").Value = 9613.429
$ws.Range("J132").Value = 8380.6
$ws.Range("K132").Value = 28840.287
$ws.Range("L132").Value = 25141.8
$ws.Range("M132").Value = -26310.287
$ws.Range("N132").Value = -30201.8
